$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextCell "D2" "27.021.17"
$ws.Range("E2").Value = "  +0.36%  "
Set-TextCell "D3" "1.678.93"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  +0.12%  "
Set-TextCell "D5" "215.88"
$ws.Range("E5").Value = "  +0.12%  "
Set-TextCell "D6" "0.519"
$ws.Range("E6").Value = "  -2.56%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("E9").Value = "  +5.36%  "
Set-TextCell "D10" "0.0625"
$ws.Range("E10").Value = "  +0.87%  "
Set-TextCell "D11" "0.0890"
$ws.Range("E11").Value = "  -0.93%  "
Set-TextCell "D12" "1.915.22"
$ws.Range("E12").Value = "  +0.77%  "
Set-TextCell "D13" "1.676.18"
$ws.Range("E13").Value = "  +0.69%  "
Set-TextCell "D14" "4.11"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("E15").Value = "  +1.74%  "
Set-TextCell "D16" "66.51"
Set-TextCell "D17" "27.026.97"
$ws.Range("E17").Value = "  +0.40%  "
Set-TextCell "D18" "8.16"
$ws.Range("E18").Value = "  +1.61%  "
Set-TextCell "D19" "236.01"
$ws.Range("E19").Value = "  +0.97%  "
Set-TextCell "D22" "4.47"
$ws.Range("E22").Value = "  +1.70%  "
Set-TextCell "D23" "9.27"
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("E24").Value = "  -4.09%  "
Set-TextCell "D25" "146.58"
$ws.Range("E25").Value = "  +0.26%  "
Set-TextCell "D26" "7.27"
$ws.Range("E26").Value = "  +1.91%  "
Set-TextCell "D27" "16.43"
$ws.Range("E27").Value = "  +3.09%  "
$ws.Range("E28").Value = "  -2.10%  "
Set-TextCell "D29" "0.999"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  +0.03%  "
Set-TextCell "D32" "3.37"
$ws.Range("E32").Value = "  +0.44%  "
Set-TextCell "D33" "1.543.39"
$ws.Range("E33").Value = "  +5.67%  "
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  +4.78%  "
$ws.Range("E36").Value = "  -0.57%  "
Set-TextCell "D37" "0.590"
$ws.Range("E37").Value = "  +2.13%  "
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("E40").Value = "  +6.42%  "
$ws.Range("E41").Value = "  +0.11%  "
Set-TextCell "D42" "67.97"
$ws.Range("E42").Value = "  +3.05%  "
Set-TextCell "D43" "5.61"
$ws.Range("E43").Value = "  -2.49%  "
Set-TextCell "D44" "2.26"
$ws.Range("E44").Value = "  -0.45%  "
Set-TextCell "D45" "1.820.38"
$ws.Range("E45").Value = "  +0.70%  "
Set-TextCell "D46" "0.780"
$ws.Range("E46").Value = "  -0.21%  "
Set-TextCell "D47" "90.35"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  +2.48%  "
Set-TextCell "D50" "7.99"
$ws.Range("E50").Value = "  +5.82%  "
$ws.Range("E51").Value = "  -0.07%  "
